$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style of A2 (bold/bordered/centered) down to the new rows A6:A21
$ws.Range("A2").Copy()
$ws.Range("A6:A21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Write data values for rows 2-21, columns A, B, C
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 657414.5417426862
$ws.Range("C2").Value = 3880.172613787232
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 428430.6638002892
$ws.Range("C3").Value = 3880.172613787232
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 306508.6679032648
$ws.Range("C4").Value = 3880.172613787232
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 186886.944987875
$ws.Range("C5").Value = 3880.172613787232
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 148289.3615605516
$ws.Range("C6").Value = 3880.172613787232
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 94964.65563556696
$ws.Range("C7").Value = 0
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 91526.37934841527
$ws.Range("C8").Value = 0
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 72265.96045345723
$ws.Range("C9").Value = 0
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 58248.21590530986
$ws.Range("C10").Value = 0
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = 53107.21270695536
$ws.Range("C11").Value = 0
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = 61653.42450662095
$ws.Range("C12").Value = 0
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = 67256.42908928309
$ws.Range("C13").Value = 0
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = 49327.45515042832
$ws.Range("C14").Value = 0
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = 49709.58303631068
$ws.Range("C15").Value = 0
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = 39133.50062070278
$ws.Range("C16").Value = 0
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = 41326.3422664292
$ws.Range("C17").Value = 0
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = 35556.4281934922
$ws.Range("C18").Value = 0
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = 31362.73578052108
$ws.Range("C19").Value = 0
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = 28169.72237710975
$ws.Range("C20").Value = 0
$ws.Range("A21").Value = 19
$ws.Range("B21").Value = 37739.54224100125
$ws.Range("C21").Value = 0

Write-Output "done"
